# Build the expanded document structure, then stamp each paragraph with
# its exact target OOXML via Range.InsertXML (keeps formatting/run-splits
# byte-precise while still going through the Word object model).

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParaXml($para, [string]$innerP) {
    $para.Range.InsertXML($pkgOpen + $innerP + $pkgClose)
}

# --- Step 1: grow the body to the final paragraph count -------------------
# Original doc has 2 paragraphs:
#   1) "{{salutation}} {{first-name}}"
#   2) a lone manual line break (<w:br/>)
# Target has 12. Paragraphs 1..8 live before the break paragraph (which
# becomes paragraph 9), paragraphs 10..12 live after it.

$p1 = $d.Paragraphs(1)
for ($i = 0; $i -lt 7; $i++) {
    $p1.Range.InsertParagraphAfter()
}

$pBreak = $d.Paragraphs(9)
for ($i = 0; $i -lt 3; $i++) {
    $pBreak.Range.InsertParagraphAfter()
}

# --- Step 2: stamp each paragraph's exact content/formatting --------------

Set-ParaXml $d.Paragraphs(1) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>{{salutation}} {{first-name}}</w:t></w:r></w:p>'

Set-ParaXml $d.Paragraphs(2) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p>'

Set-ParaXml $d.Paragraphs(3) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>{{salutation}} {{first-name}}</w:t></w:r></w:p>'

Set-ParaXml $d.Paragraphs(4) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p>'

Set-ParaXml $d.Paragraphs(5) '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Test</w:t></w:r></w:p>'

Set-ParaXml $d.Paragraphs(6) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p>'

Set-ParaXml $d.Paragraphs(7) '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr><w:t>{{salutation}} {{first-name}}</w:t></w:r></w:p>'

Set-ParaXml $d.Paragraphs(8) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p>'

Set-ParaXml $d.Paragraphs(9) '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:sectPr><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1417" w:right="1417" w:bottom="1134" w:left="1417" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/></w:sectPr></w:pPr><w:r><w:br/></w:r></w:p>'

Set-ParaXml $d.Paragraphs(10) '<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:color w:val="C0504D" w:themeColor="accent2"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="C0504D" w:themeColor="accent2"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:lastRenderedPageBreak/><w:t>{{salutation}} {{first-name}}</w:t></w:r></w:p>'

Set-ParaXml $d.Paragraphs(11) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p>'

Set-ParaXml $d.Paragraphs(12) '<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>The person {{salut</w:t></w:r><w:r><w:t>at</w:t></w:r><w:r><w:t>ion}} is a {{salut</w:t></w:r><w:r><w:t>at</w:t></w:r><w:r><w:t>ion}}</w:t></w:r></w:p>'

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
